# "fixed new assigned team to user"
# Updates the BU/Contractor Code/Planner Group/Planner Center Name/Contractor
# assignment on the "Create Teams" sheet, and adds a new assignment row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create Teams")

# Correct the existing assignment in row 2 (was 0-BR-MVS-02 / EPVL0001 / ZP1 / 3R71 / VOLTALIA)
$ws.Range("A2").Value = "0-ES-BEL-01"
$ws.Range("B2").Value = "ETVT0001"
$ws.Range("C2").Value = "ZP1"
$ws.Range("D2").Value = "H314"
$ws.Range("E2").Value = "VESTAS"

# Add the new assignment row 3
$ws.Range("A3").Value = "0-ES-ACA-01"
$ws.Range("B3").Value = "ETVT0001"
$ws.Range("C3").Value = "ZP1"
$ws.Range("D3").Value = "T558"
$ws.Range("E3").Value = "VESTAS"

# Leave the selection where the author ended up after editing
$ws.Range("C13").Select() | Out-Null
